$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The commit adds one new daily-price record for "Apio" (Vega Modelo de
# Temuco) right before the existing row 294, pushing every subsequent
# record down by one row (old row 392 becomes new row 393).
$ws.Rows("294").Insert()

$row = 294
$ws.Cells.Item($row,1).Value  = 10
$ws.Cells.Item($row,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item($row,3).Value  = "La Araucanía"
$ws.Cells.Item($row,4).Value  = 44876
$ws.Cells.Item($row,5).Value  = 9
$ws.Cells.Item($row,6).Value  = 100112017
$ws.Cells.Item($row,7).Value  = "Apio"
$ws.Cells.Item($row,8).Value  = "Americana (o)"
$ws.Cells.Item($row,9).Value  = "Primera"
$ws.Cells.Item($row,10).Value = 95
$ws.Cells.Item($row,11).Value = 9000
$ws.Cells.Item($row,12).Value = 10000
$ws.Cells.Item($row,13).Value = 9421
$ws.Cells.Item($row,14).Value = "$/docena de matas"
$ws.Cells.Item($row,15).Value = "Provincia del Elquí"
$ws.Cells.Item($row,16).Value = 1570
$ws.Cells.Item($row,17).Value = 6
$ws.Cells.Item($row,18).Value = "Hortaliza"
